# Apply numeric cell updates produced by the scheduled profit-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1345.8518
$ws.Range("J17").Value = 1345.8518
$ws.Range("L17").Value = 4037.5554
$ws.Range("N17").Value = -4373.555399999999
$ws.Range("H43").Value = 458338.22
$ws.Range("I43").Value = 498.33334
$ws.Range("K43").Value = 498.33334
$ws.Range("M43").Value = -429.33334
$ws.Range("H80").Value = 9654475
$ws.Range("I80").Value = 20850212
$ws.Range("J80").Value = 58129
$ws.Range("K80").Value = 62550636
$ws.Range("L80").Value = 174387
$ws.Range("M80").Value = -62549638
$ws.Range("N80").Value = -176383
$ws.Range("H83").Value = 9654475
$ws.Range("I83").Value = 20850212
$ws.Range("J83").Value = 58129
$ws.Range("K83").Value = 187651908
$ws.Range("L83").Value = 523161
$ws.Range("M83").Value = -187646916
$ws.Range("N83").Value = -533145
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H116").Value = 13167764
$ws.Range("I116").Value = 27784058
$ws.Range("J116").Value = 13098.9
$ws.Range("K116").Value = 27784058
$ws.Range("L116").Value = 13098.9
$ws.Range("M116").Value = -27780616
$ws.Range("N116").Value = -19982.9
$ws.Range("H132").Value = 1083.7142
$ws.Range("I132").Value = 1027.9
$ws.Range("K132").Value = 3083.7
$ws.Range("M132").Value = -553.7000000000003
$ws.Range("H137").Value = 2942.6296
$ws.Range("J137").Value = 4170.7144
$ws.Range("L137").Value = 12512.1432
$ws.Range("N137").Value = -17612.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4088349.2
$ws.Range("I32").Value = 4449336
$ws.Range("K32").Value = 4449336
$ws.Range("M32").Value = -4449049
$ws.Range("H45").Value = 1923.2609
$ws.Range("I45").Value = 1922.6154
$ws.Range("K45").Value = 1922.6154
$ws.Range("M45").Value = -1545.6154
$ws.Range("H61").Value = 41669956
$ws.Range("I61").Value = 1470.2142
$ws.Range("K61").Value = 1470.2142
$ws.Range("M61").Value = -1258.2142
$ws.Range("H74").Value = 32935.727
$ws.Range("I74").Value = 41634.08
$ws.Range("J74").Value = 5753.375
$ws.Range("K74").Value = 41634.08
$ws.Range("L74").Value = 5753.375
$ws.Range("M74").Value = -40760.08
$ws.Range("N74").Value = -7501.375
$ws.Range("H77").Value = 32935.727
$ws.Range("I77").Value = 41634.08
$ws.Range("J77").Value = 5753.375
$ws.Range("K77").Value = 208170.4
$ws.Range("L77").Value = 28766.875
$ws.Range("M77").Value = -203802.4
$ws.Range("N77").Value = -37502.875
$ws.Range("H110").Value = 66668030
$ws.Range("J110").Value = 66668030
$ws.Range("L110").Value = 66668030
$ws.Range("N110").Value = -66672120
$ws.Range("H122").Value = 4408.3335
$ws.Range("J122").Value = 7037.273
$ws.Range("L122").Value = 21111.819
$ws.Range("N122").Value = -26011.819
$ws.Range("H132").Value = 4695.094
$ws.Range("I132").Value = 3866.8333
$ws.Range("J132").Value = 6449.0586
$ws.Range("K132").Value = 11600.4999
$ws.Range("L132").Value = 19347.1758
$ws.Range("M132").Value = -9070.499899999999
$ws.Range("N132").Value = -24407.1758
$ws.Range("H136").Value = 41669956
$ws.Range("I136").Value = 1470.2142
$ws.Range("K136").Value = 4410.642599999999
$ws.Range("M136").Value = -1860.642599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I134").Value = 7813918
$ws.Range("K134").Value = 23441754
$ws.Range("M134").Value = -23439219

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4655.525
$ws.Range("J31").Value = 6896
$ws.Range("L31").Value = 6896
$ws.Range("N31").Value = -7486
$ws.Range("H34").Value = 4655.525
$ws.Range("J34").Value = 6896
$ws.Range("L34").Value = 6896
$ws.Range("N34").Value = -7300
$ws.Range("H62").Value = 4637.8076
$ws.Range("I62").Value = 4376.2354
$ws.Range("J62").Value = 5131.8887
$ws.Range("K62").Value = 4376.2354
$ws.Range("L62").Value = 5131.8887
$ws.Range("M62").Value = -3752.2354
$ws.Range("N62").Value = -6379.8887
$ws.Range("H65").Value = 4637.8076
$ws.Range("I65").Value = 4376.2354
$ws.Range("J65").Value = 5131.8887
$ws.Range("K65").Value = 21881.177
$ws.Range("L65").Value = 25659.4435
$ws.Range("M65").Value = -18761.177
$ws.Range("N65").Value = -31899.4435
$ws.Range("H107").Value = 3286
$ws.Range("I107").Value = 3431.3333
$ws.Range("J107").Value = 3223.7144
$ws.Range("K107").Value = 3431.3333
$ws.Range("L107").Value = 3223.7144
$ws.Range("M107").Value = -1511.3333
$ws.Range("N107").Value = -7063.7144
$ws.Range("H122").Value = 4698.2744
$ws.Range("I122").Value = 4158.7715
$ws.Range("J122").Value = 5878.4375
$ws.Range("K122").Value = 12476.3145
$ws.Range("L122").Value = 17635.3125
$ws.Range("M122").Value = -10026.3145
$ws.Range("N122").Value = -22535.3125
$ws.Range("H132").Value = 3545.5122
$ws.Range("I132").Value = 2729.6296
$ws.Range("K132").Value = 8188.888800000001
$ws.Range("M132").Value = -5658.888800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 704192.25
$ws.Range("I4").Value = 2460.375
$ws.Range("J4").Value = 1639834.8
$ws.Range("K4").Value = 7381.125
$ws.Range("L4").Value = 4919504.4
$ws.Range("M4").Value = -7269.125
$ws.Range("N4").Value = -4919728.4
$ws.Range("H131").Value = 2565.6042
$ws.Range("J131").Value = 2535.432
$ws.Range("L131").Value = 7606.295999999999
$ws.Range("N131").Value = -17686.296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H127").Value = 42000
$ws.Range("J127").Value = 42000
$ws.Range("L127").Value = 42000
$ws.Range("N127").Value = -51920
$ws.Range("H128").Value = 59500
$ws.Range("J128").Value = 59500
$ws.Range("L128").Value = 59500
$ws.Range("N128").Value = -69460
$ws.Range("H130").Value = 74333.336
$ws.Range("I130").Value = 40000
$ws.Range("J130").Value = 91500
$ws.Range("K130").Value = 40000
$ws.Range("L130").Value = 91500
$ws.Range("M130").Value = -34980
$ws.Range("N130").Value = -101540
$ws.Range("H132").Value = 3229.923
$ws.Range("I132").Value = 2144.375
$ws.Range("K132").Value = 6433.125
$ws.Range("M132").Value = -3903.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 396
$ws.Range("I16").Value = 411.66666
$ws.Range("J16").Value = 302
$ws.Range("K16").Value = 411.66666
$ws.Range("L16").Value = 302
$ws.Range("M16").Value = -241.66666
$ws.Range("N16").Value = -642
$ws.Range("H40").Value = 35718650
$ws.Range("I40").Value = 55557900
$ws.Range("K40").Value = 55557900
$ws.Range("M40").Value = -55557764
$ws.Range("H60").Value = 34557.668
$ws.Range("J60").Value = 34557.668
$ws.Range("L60").Value = 34557.668
$ws.Range("N60").Value = -35575.668
$ws.Range("H61").Value = 5611.476
$ws.Range("I61").Value = 4505.636
$ws.Range("J61").Value = 6827.9
$ws.Range("K61").Value = 4505.636
$ws.Range("L61").Value = 6827.9
$ws.Range("M61").Value = -4303.636
$ws.Range("N61").Value = -7231.9
$ws.Range("H113").Value = 5611.476
$ws.Range("I113").Value = 4505.636
$ws.Range("J113").Value = 6827.9
$ws.Range("K113").Value = 4505.636
$ws.Range("L113").Value = 6827.9
$ws.Range("M113").Value = -2335.636
$ws.Range("N113").Value = -11167.9
$ws.Range("H122").Value = 4262.4165
$ws.Range("J122").Value = 5145.5
$ws.Range("L122").Value = 15436.5
$ws.Range("N122").Value = -20336.5
$ws.Range("H132").Value = 8071699
$ws.Range("I132").Value = 15154528
$ws.Range("K132").Value = 45463584
$ws.Range("M132").Value = -45461054
$ws.Range("H136").Value = 7765.8228
$ws.Range("I136").Value = 2683.2856
$ws.Range("K136").Value = 8049.8568
$ws.Range("M136").Value = -5499.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4708
$ws.Range("H107").Value = 27779498
$ws.Range("I107").Value = 1153.3334
$ws.Range("J107").Value = 55557840
$ws.Range("K107").Value = 3460.0002
$ws.Range("L107").Value = 166673520
$ws.Range("M107").Value = -1540.0002
$ws.Range("N107").Value = -166677360
$ws.Range("H122").Value = 20165032
$ws.Range("I122").Value = 29650260
$ws.Range("K122").Value = 88950780
$ws.Range("M122").Value = -88948330
$ws.Range("H132").Value = 4272.6
$ws.Range("I132").Value = 4074.9688
$ws.Range("J132").Value = 5063.125
$ws.Range("K132").Value = 12224.9064
$ws.Range("L132").Value = 15189.375
$ws.Range("M132").Value = -9694.9064
$ws.Range("N132").Value = -20249.375
$ws.Range("H136").Value = 23491774
$ws.Range("I136").Value = 47619880
$ws.Range("K136").Value = 142859640
$ws.Range("M136").Value = -142857090
